$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 950
$ws.Range("I12").Value = 950
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 950
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -780
$ws.Range("H75").Value = 29896.75
$ws.Range("J75").Value = 29896.75
$ws.Range("L75").Value = 29896.75
$ws.Range("N75").Value = -31768.75
$ws.Range("H78").Value = 29896.75
$ws.Range("J78").Value = 29896.75
$ws.Range("L78").Value = 89690.25
$ws.Range("N78").Value = -99050.25
$ws.Range("H116").Value = 8419.786
$ws.Range("I116").Value = 8222.5
$ws.Range("K116").Value = 8222.5
$ws.Range("M116").Value = -4780.5
$ws.Range("H125").Value = 2601
$ws.Range("J125").Value = 2971.111
$ws.Range("L125").Value = 26739.999
$ws.Range("N125").Value = -31659.999
$ws.Range("H135").Value = 2000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070
$ws.Range("H137").Value = 1547.7368
$ws.Range("I137").Value = 1245.9
$ws.Range("K137").Value = 3737.7
$ws.Range("M137").Value = -1187.7
$ws.Range("H138").Value = 6976.7856
$ws.Range("I138").Value = 7517.6523
$ws.Range("K138").Value = 22552.9569
$ws.Range("M138").Value = -17412.9569
$ws.Range("N12").ClearContents()
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2337.4119
$ws.Range("I2").Value = 976.9
$ws.Range("J2").Value = 4281
$ws.Range("K2").Value = 976.9
$ws.Range("L2").Value = 4281
$ws.Range("M2").Value = -863.9
$ws.Range("N2").Value = -4507
$ws.Range("H43").Value = 19377
$ws.Range("J43").Value = 19377
$ws.Range("L43").Value = 19377
$ws.Range("N43").Value = -20003
$ws.Range("H45").Value = 72982.92999999999
$ws.Range("I45").Value = 92414.63
$ws.Range("K45").Value = 92414.63
$ws.Range("M45").Value = -92037.63
$ws.Range("H63").Value = 13615.728
$ws.Range("I63").Value = 7795.6665
$ws.Range("J63").Value = 20599.8
$ws.Range("K63").Value = 7795.6665
$ws.Range("L63").Value = 20599.8
$ws.Range("M63").Value = -7109.6665
$ws.Range("N63").Value = -21971.8
$ws.Range("H66").Value = 13615.728
$ws.Range("I66").Value = 7795.6665
$ws.Range("J66").Value = 20599.8
$ws.Range("K66").Value = 38978.3325
$ws.Range("L66").Value = 102999
$ws.Range("M66").Value = -35546.3325
$ws.Range("N66").Value = -109863
$ws.Range("H102").Value = 4281.5
$ws.Range("I102").Value = 4783.6665
$ws.Range("K102").Value = 4783.6665
$ws.Range("M102").Value = -3161.6665
$ws.Range("H116").Value = 2337.4119
$ws.Range("I116").Value = 976.9
$ws.Range("J116").Value = 4281
$ws.Range("K116").Value = 976.9
$ws.Range("L116").Value = 4281
$ws.Range("M116").Value = 1317.1
$ws.Range("N116").Value = -8869

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2337.4119
$ws.Range("I3").Value = 976.9
$ws.Range("J3").Value = 4281
$ws.Range("K3").Value = 976.9
$ws.Range("L3").Value = 4281
$ws.Range("M3").Value = -862.9
$ws.Range("N3").Value = -4509
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H82").Value = 17686.111
$ws.Range("I82").Value = 11771.875
$ws.Range("J82").Value = 65000
$ws.Range("K82").Value = 11771.875
$ws.Range("L82").Value = 65000
$ws.Range("M82").Value = -11388.875
$ws.Range("N82").Value = -65766
$ws.Range("H85").Value = 17686.111
$ws.Range("I85").Value = 11771.875
$ws.Range("J85").Value = 65000
$ws.Range("K85").Value = 11771.875
$ws.Range("L85").Value = 65000
$ws.Range("M85").Value = -10445.875
$ws.Range("N85").Value = -67652
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6911.3125
$ws.Range("I31").Value = 4729.5
$ws.Range("K31").Value = 4729.5
$ws.Range("M31").Value = -4434.5
$ws.Range("H34").Value = 6911.3125
$ws.Range("I34").Value = 4729.5
$ws.Range("K34").Value = 4729.5
$ws.Range("M34").Value = -4527.5
$ws.Range("H58").Value = 2352.8333
$ws.Range("I58").Value = 2324.4
$ws.Range("K58").Value = 2324.4
$ws.Range("M58").Value = -2121.4
$ws.Range("H88").Value = 8198.5
$ws.Range("J88").Value = 8198.5
$ws.Range("L88").Value = 8198.5
$ws.Range("N88").Value = -9010.5
$ws.Range("H91").Value = 8198.5
$ws.Range("J91").Value = 8198.5
$ws.Range("L91").Value = 8198.5
$ws.Range("N91").Value = -11006.5
$ws.Range("H132").Value = 4172.773
$ws.Range("I132").Value = 4659.5
$ws.Range("J132").Value = 3767.1667
$ws.Range("K132").Value = 13978.5
$ws.Range("L132").Value = 11301.5001
$ws.Range("M132").Value = -11448.5
$ws.Range("N132").Value = -16361.5001
$ws.Range("H134").Value = 2573
$ws.Range("I134").Value = 2355
$ws.Range("K134").Value = 7065
$ws.Range("M134").Value = -4530
$ws.Range("H136").Value = 2352.8333
$ws.Range("I136").Value = 2324.4
$ws.Range("K136").Value = 6973.200000000001
$ws.Range("M136").Value = -4423.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1281349.6
$ws.Range("I4").Value = 1824973
$ws.Range("K4").Value = 5474919
$ws.Range("M4").Value = -5474807
$ws.Range("H23").Value = 100126.2
$ws.Range("I23").Value = 143.5
$ws.Range("K23").Value = 430.5
$ws.Range("M23").Value = -195.5
$ws.Range("H60").Value = 743.2
$ws.Range("I60").Value = 774.8570999999999
$ws.Range("K60").Value = 2324.5713
$ws.Range("M60").Value = -2073.5713
$ws.Range("H104").Value = 7795.162
$ws.Range("J104").Value = 7811.7427
$ws.Range("L104").Value = 23435.2281
$ws.Range("N104").Value = -28677.2281

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 45487.5
$ws.Range("I62").Value = 38000
$ws.Range("K62").Value = 38000
$ws.Range("M62").Value = -37314
$ws.Range("H65").Value = 45487.5
$ws.Range("I65").Value = 38000
$ws.Range("K65").Value = 114000
$ws.Range("M65").Value = -110568
$ws.Range("H80").Value = 129950.625
$ws.Range("J80").Value = 42904.94
$ws.Range("L80").Value = 42904.94
$ws.Range("N80").Value = -44900.94
$ws.Range("H83").Value = 129950.625
$ws.Range("J83").Value = 42904.94
$ws.Range("L83").Value = 214524.7
$ws.Range("N83").Value = -224508.7
$ws.Range("H92").Value = 12177.333
$ws.Range("J92").Value = 12177.333
$ws.Range("L92").Value = 12177.333
$ws.Range("N92").Value = -15921.333
$ws.Range("H97").Value = 81968.69500000001
$ws.Range("I97").Value = 63902.438
$ws.Range("K97").Value = 63902.438
$ws.Range("M97").Value = -63406.438

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9105.08
$ws.Range("I7").Value = 4912
$ws.Range("J7").Value = 13647.583
$ws.Range("K7").Value = 4912
$ws.Range("L7").Value = 13647.583
$ws.Range("M7").Value = -4800
$ws.Range("N7").Value = -13871.583
$ws.Range("H42").Value = 74583.336
$ws.Range("I42").Value = 82500
$ws.Range("K42").Value = 82500
$ws.Range("M42").Value = -81937
$ws.Range("H49").Value = 74583.336
$ws.Range("I49").Value = 82500
$ws.Range("K49").Value = 82500
$ws.Range("M49").Value = -82353
$ws.Range("H58").Value = 1975
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H93").Value = 3256.25
$ws.Range("I93").Value = 1210
$ws.Range("K93").Value = 1210
$ws.Range("M93").Value = 38
$ws.Range("H94").Value = 23800.572
$ws.Range("I94").Value = 20444
$ws.Range("J94").Value = 24360
$ws.Range("K94").Value = 20444
$ws.Range("L94").Value = 24360
$ws.Range("M94").Value = -19768
$ws.Range("N94").Value = -25712
$ws.Range("H100").Value = 5031.4
$ws.Range("I100").Value = 3385.6667
$ws.Range("K100").Value = 3385.6667
$ws.Range("M100").Value = -2844.6667
$ws.Range("H126").Value = 9105.08
$ws.Range("I126").Value = 4912
$ws.Range("J126").Value = 13647.583
$ws.Range("K126").Value = 14736
$ws.Range("L126").Value = 40942.749
$ws.Range("M126").Value = -12266
$ws.Range("N126").Value = -45882.749
$ws.Range("H132").Value = 8294.925999999999
$ws.Range("I132").Value = 2592.4167
$ws.Range("K132").Value = 7777.250100000001
$ws.Range("M132").Value = -5247.250100000001
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2664.7083
$ws.Range("I126").Value = 2989.1538
$ws.Range("J126").Value = 2281.2727
$ws.Range("K126").Value = 8967.4614
$ws.Range("L126").Value = 6843.8181
$ws.Range("M126").Value = -6497.4614
$ws.Range("N126").Value = -11783.8181
$ws.Range("H132").Value = 26217.477
$ws.Range("I132").Value = 33222.062
$ws.Range("K132").Value = 99666.18599999999
$ws.Range("M132").Value = -97136.18599999999
$ws.Range("H140").Value = 85724.664
$ws.Range("J140").Value = 85724.664
$ws.Range("L140").Value = 85724.664
$ws.Range("N140").Value = -96084.664
